$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("15-11-2021 11:18", "Xexe"),
    @("15-11-2021 11:18", "Fv"),
    @("15-11-2021 12:00", "B"),
    @("15-11-2021 12:51", "C"),
    @("15-11-2021 01:47", "Xbxjd")
)

$startRow = 14
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}
